$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I2").Value = "2000|0001|Y|R|0|1|0|4|2|2|2"
$ws.Rows.Item(2).RowHeight = 105
$ws.Range("I2").Select() | Out-Null
